# Updated symbol list (refresh of price/volume/hour columns, plus a
# one-row cyclic shuffle of the rank-8..18 coin block). All written
# values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the sheet's existing inlineStr/text cell type)
# instead of auto-coercing numeric- or percent-looking strings into
# real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'247.00"
$ws.Cells.Item(2, 5).Value = "'0.70%"
$ws.Cells.Item(2, 7).Value = "'22"
$ws.Cells.Item(3, 4).Value = "'26.23"
$ws.Cells.Item(3, 5).Value = "'4.37%"
$ws.Cells.Item(3, 7).Value = "'22"
$ws.Cells.Item(4, 4).Value = "'5.087"
$ws.Cells.Item(4, 5).Value = "'1.00%"
$ws.Cells.Item(4, 7).Value = "'22"
$ws.Cells.Item(5, 4).Value = "'0.05599"
$ws.Cells.Item(5, 5).Value = "'-0.23%"
$ws.Cells.Item(5, 7).Value = "'22"
$ws.Cells.Item(6, 4).Value = "'6.478"
$ws.Cells.Item(6, 5).Value = "'-1.21%"
$ws.Cells.Item(6, 7).Value = "'22"
$ws.Cells.Item(7, 4).Value = "'0.8130"
$ws.Cells.Item(7, 5).Value = "'-0.15%"
$ws.Cells.Item(7, 7).Value = "'22"
$ws.Cells.Item(8, 4).Value = "'0.8447"
$ws.Cells.Item(8, 5).Value = "'0.86%"
$ws.Cells.Item(8, 7).Value = "'22"
$ws.Cells.Item(9, 2).Value = "'MandalaExchangeToken"
$ws.Cells.Item(9, 3).Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(9, 4).Value = "'0.06988"
$ws.Cells.Item(9, 5).Value = "'0.65%"
$ws.Cells.Item(9, 7).Value = "'22"
$ws.Cells.Item(10, 2).Value = "'BitrueCoin"
$ws.Cells.Item(10, 3).Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(10, 4).Value = "'0.02813"
$ws.Cells.Item(10, 5).Value = "'-0.88%"
$ws.Cells.Item(10, 7).Value = "'22"
$ws.Cells.Item(11, 2).Value = "'BitMartToken"
$ws.Cells.Item(11, 3).Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(11, 4).Value = "'0.09393"
$ws.Cells.Item(11, 5).Value = "'-0.14%"
$ws.Cells.Item(11, 7).Value = "'22"
$ws.Cells.Item(12, 2).Value = "'BitForexToken"
$ws.Cells.Item(12, 3).Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(12, 4).Value = "'0.001513"
$ws.Cells.Item(12, 5).Value = "'0.21%"
$ws.Cells.Item(12, 7).Value = "'22"
$ws.Cells.Item(13, 2).Value = "'One"
$ws.Cells.Item(13, 3).Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(13, 4).Value = "'0.0005992"
$ws.Cells.Item(13, 5).Value = "'0.23%"
$ws.Cells.Item(13, 7).Value = "'22"
$ws.Cells.Item(14, 2).Value = "'TigerCash"
$ws.Cells.Item(14, 3).Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(14, 4).Value = "'0.006163"
$ws.Cells.Item(14, 5).Value = "'-0.04%"
$ws.Cells.Item(14, 7).Value = "'22"
$ws.Cells.Item(15, 2).Value = "'LEO"
$ws.Cells.Item(15, 3).Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(15, 4).Value = "'3.608"
$ws.Cells.Item(15, 5).Value = "'3.11%"
$ws.Cells.Item(15, 7).Value = "'22"
$ws.Cells.Item(16, 2).Value = "'GateToken"
$ws.Cells.Item(16, 3).Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(16, 4).Value = "'3.018"
$ws.Cells.Item(16, 5).Value = "'0.25%"
$ws.Cells.Item(16, 7).Value = "'22"
$ws.Cells.Item(17, 2).Value = "'BTSEToken"
$ws.Cells.Item(17, 3).Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(17, 4).Value = "'2.055"
$ws.Cells.Item(17, 5).Value = "'-1.74%"
$ws.Cells.Item(17, 7).Value = "'22"
$ws.Cells.Item(18, 2).Value = "'BitpandaEcosystemToken"
$ws.Cells.Item(18, 3).Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(18, 4).Value = "'0.3113"
$ws.Cells.Item(18, 5).Value = "'-1.48%"
$ws.Cells.Item(18, 7).Value = "'22"
$ws.Cells.Item(19, 2).Value = "'WazirX"
$ws.Cells.Item(19, 3).Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(19, 4).Value = "'0.1336"
$ws.Cells.Item(19, 5).Value = "'-0.22%"
$ws.Cells.Item(19, 7).Value = "'22"
$ws.Cells.Item(20, 4).Value = "'0.03207"
$ws.Cells.Item(20, 5).Value = "'-1.72%"
$ws.Cells.Item(20, 7).Value = "'22"
$ws.Cells.Item(21, 5).Value = "'-1.28%"
$ws.Cells.Item(21, 7).Value = "'22"
$ws.Cells.Item(22, 4).Value = "'3.753"
$ws.Cells.Item(22, 5).Value = "'0.11%"
$ws.Cells.Item(22, 7).Value = "'22"
$ws.Cells.Item(23, 4).Value = "'0.04644"
$ws.Cells.Item(23, 5).Value = "'-1.02%"
$ws.Cells.Item(23, 7).Value = "'22"
$ws.Cells.Item(24, 5).Value = "'0.40%"
$ws.Cells.Item(24, 7).Value = "'22"
$ws.Cells.Item(25, 4).Value = "'0.001246"
$ws.Cells.Item(25, 5).Value = "'0.40%"
$ws.Cells.Item(25, 7).Value = "'22"
$ws.Cells.Item(26, 4).Value = "'0.004575"
$ws.Cells.Item(26, 5).Value = "'1.05%"
$ws.Cells.Item(26, 7).Value = "'22"
$ws.Cells.Item(27, 4).Value = "'0.00009602"
$ws.Cells.Item(27, 5).Value = "'-1.02%"
$ws.Cells.Item(27, 7).Value = "'22"
$ws.Cells.Item(28, 4).Value = "'0.0001939"
$ws.Cells.Item(28, 5).Value = "'-0.03%"
$ws.Cells.Item(28, 7).Value = "'22"
$ws.Cells.Item(29, 7).Value = "'22"
$ws.Cells.Item(30, 7).Value = "'22"
$ws.Cells.Item(31, 7).Value = "'22"
$ws.Cells.Item(32, 7).Value = "'22"
$ws.Cells.Item(33, 7).Value = "'22"
$ws.Cells.Item(34, 7).Value = "'22"
$ws.Cells.Item(35, 7).Value = "'22"
$ws.Cells.Item(36, 7).Value = "'22"
$ws.Cells.Item(37, 7).Value = "'22"
$ws.Cells.Item(38, 7).Value = "'22"
$ws.Cells.Item(39, 7).Value = "'22"
$ws.Cells.Item(40, 4).Value = "'0.03662"
$ws.Cells.Item(40, 5).Value = "'-0.08%"
$ws.Cells.Item(40, 7).Value = "'22"
$ws.Cells.Item(41, 4).Value = "'0.003427"
$ws.Cells.Item(41, 5).Value = "'-45.08%"
$ws.Cells.Item(41, 7).Value = "'22"
$ws.Cells.Item(42, 4).Value = "'0.1355"
$ws.Cells.Item(42, 5).Value = "'-0.09%"
$ws.Cells.Item(42, 7).Value = "'22"
$ws.Cells.Item(43, 4).Value = "'0.002660"
$ws.Cells.Item(43, 5).Value = "'-2.72%"
$ws.Cells.Item(43, 7).Value = "'22"
$ws.Cells.Item(44, 4).Value = "'0.008261"
$ws.Cells.Item(44, 5).Value = "'0.79%"
$ws.Cells.Item(44, 7).Value = "'22"
$ws.Cells.Item(45, 4).Value = "'0.00005394"
$ws.Cells.Item(45, 5).Value = "'1.87%"
$ws.Cells.Item(45, 7).Value = "'22"
$ws.Cells.Item(46, 5).Value = "'0.07%"
$ws.Cells.Item(46, 7).Value = "'22"
$ws.Cells.Item(47, 5).Value = "'-35.80%"
$ws.Cells.Item(47, 7).Value = "'22"
$ws.Cells.Item(48, 5).Value = "'20.58%"
$ws.Cells.Item(48, 7).Value = "'22"
$ws.Cells.Item(49, 4).Value = "'0.00002101"
$ws.Cells.Item(49, 5).Value = "'0.07%"
$ws.Cells.Item(49, 7).Value = "'22"
$ws.Cells.Item(50, 4).Value = "'0.0002001"
$ws.Cells.Item(50, 5).Value = "'0.07%"
$ws.Cells.Item(50, 7).Value = "'22"
$ws.Cells.Item(51, 7).Value = "'22"
